# Update "想去人数" (want-to-go count) figures in column F for the two
# sheets that carry exhibition data ("展览" and "全部类型"). Both sheets
# mirror the same rows, so the same F-column edits apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 609
    6  = 1092
    7  = 9
    8  = 11372
    12 = 348
    14 = 785
    15 = 12338
    16 = 13005
    17 = 33
    23 = 94
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F" + $row).Value = $updates[$row]
    }
}
